$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row below the existing blank spacer row (row 9), pushing the
# totals block (old rows 9-12) down by one row. This keeps the SUM formula's
# range (C4:C9) textually the same while the new data row becomes row 9.
$ws.Rows(10).Insert()

# Copy formatting from the row above (row 8) into the freshly inserted row 9
# so the new cells pick up the same styles (date format, centered, etc.)
$ws.Range("A8:D8").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)

# Fill in the new timeline entry (Day 6)
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = 45570
$ws.Range("C9").Value = 1.5
$ws.Range("D9").Value = "Added Username reset and integrated with new phases functionality on BE"

# Correct the previously wrong dates for Day 4 and Day 5
$ws.Range("B7").Value = 45296
$ws.Range("B8").Value = 45509

# Update the saved selection shown in the sheet view
$ws.Range("D17").Select()
